# Apply "repull data, push all data, mean calculation" edit:
# Update the dSF column (F) values for specific rows to reflect
# recalculated/repulled data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    2  = -5
    7  = 2
    8  = 0
    13 = 0
    16 = 1
    17 = 1
    18 = 3
    33 = -1
    43 = 1
    49 = 0
    53 = 2
    54 = -1
    55 = 2
    57 = 3
    58 = 0
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
